$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 3,5,7,9,11,13,15,17,19,21) {
    $ws.Range("D$r").Value = 4
}
